# Generate Report for Handback
#
# This script applies a "handback" update to the localization-status workbook:
#   - Sheet "Overview": Status columns (zh-cn / de-de) change from
#     "Ready for handoff" to "Handed back: in sync with en-US"
#   - Sheet "zh-cn": Status changes the same way, a "Latest Target File" (F)
#     and "Latest Handback File" (G) hyperlink are filled in for both data
#     rows, and "Latest Handback DateTime" (H) is stamped.
#   - Sheet "de-de": same treatment, with its own handback timestamp.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# Color used by the workbook's existing hyperlink cells (RGB 0x6495ED,
# encoded for the Font.Color COM property as 0x00BBGGRR).
$hyperlinkColor = 15570276

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Name = "Calibri"
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aefd90d31dcddd98c512ca71578c6135ada0f601/e2e/2ca97931-b6db-43bf-a64b-b17f3dd3d1f8.md",
    "",
    "",
    "2ca97931-b6db-43bf-a64b-b17f3dd3d1f8.md")
Style-AsHyperlink $wsZhCn.Range("F2")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7aefdda2fdc70b79906983a2cb336d6bcfe45f28/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/2ca97931-b6db-43bf-a64b-b17f3dd3d1f8.60940da3589336e71e27ad785bd9d9f888631f00.zh-cn.xlf",
    "",
    "",
    "2ca97931-b6db-43bf-a64b-b17f3dd3d1f8.60940da3589336e71e27ad785bd9d9f888631f00.zh-cn.xlf")
Style-AsHyperlink $wsZhCn.Range("G2")

$wsZhCn.Range("H2").Value = "2016-03-22 05:50:28"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aefd90d31dcddd98c512ca71578c6135ada0f601/e2e/30ea4186-a07d-49a5-a220-82dc7c32336c.md",
    "",
    "",
    "30ea4186-a07d-49a5-a220-82dc7c32336c.md")
Style-AsHyperlink $wsZhCn.Range("F3")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7aefdda2fdc70b79906983a2cb336d6bcfe45f28/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/30ea4186-a07d-49a5-a220-82dc7c32336c.791958e343af07176d13bc030ea75dd550d998e5.zh-cn.xlf",
    "",
    "",
    "30ea4186-a07d-49a5-a220-82dc7c32336c.791958e343af07176d13bc030ea75dd550d998e5.zh-cn.xlf")
Style-AsHyperlink $wsZhCn.Range("G3")

$wsZhCn.Range("H3").Value = "2016-03-22 05:50:28"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aefd90d31dcddd98c512ca71578c6135ada0f601/e2e/2ca97931-b6db-43bf-a64b-b17f3dd3d1f8.md",
    "",
    "",
    "2ca97931-b6db-43bf-a64b-b17f3dd3d1f8.md")
Style-AsHyperlink $wsDeDe.Range("F2")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d387566ccd8d4a81adc128580ae8cf5463f57d8b/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/2ca97931-b6db-43bf-a64b-b17f3dd3d1f8.60940da3589336e71e27ad785bd9d9f888631f00.de-de.xlf",
    "",
    "",
    "2ca97931-b6db-43bf-a64b-b17f3dd3d1f8.60940da3589336e71e27ad785bd9d9f888631f00.de-de.xlf")
Style-AsHyperlink $wsDeDe.Range("G2")

$wsDeDe.Range("H2").Value = "2016-03-22 05:50:42"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/aefd90d31dcddd98c512ca71578c6135ada0f601/e2e/30ea4186-a07d-49a5-a220-82dc7c32336c.md",
    "",
    "",
    "30ea4186-a07d-49a5-a220-82dc7c32336c.md")
Style-AsHyperlink $wsDeDe.Range("F3")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d387566ccd8d4a81adc128580ae8cf5463f57d8b/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/30ea4186-a07d-49a5-a220-82dc7c32336c.791958e343af07176d13bc030ea75dd550d998e5.de-de.xlf",
    "",
    "",
    "30ea4186-a07d-49a5-a220-82dc7c32336c.791958e343af07176d13bc030ea75dd550d998e5.de-de.xlf")
Style-AsHyperlink $wsDeDe.Range("G3")

$wsDeDe.Range("H3").Value = "2016-03-22 05:50:42"
